# "Oprydning og samle i 'Spil hele lortet'" - update the Rounded_CH (E)
# and the unnamed 5th (F) columns with the recomputed Nash-equilibrium
# values for each deck (rows 2-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ E = 1.596407552675339; F = 1.596407552675339 }
    3  = @{ E = 1.596407552675339; F = 1.596407552675339 }
    4  = @{ E = 1.596407552675339; F = 1.596407552675339 }
    5  = @{ E = 22.44744497537364; F = 1.596407552675339 }
    6  = @{ E = 38.08572304239735; F = 19.84106529753635 }
    7  = @{ E = 1.596407552675339; F = 1.596407552675339 }
    8  = @{ E = 1.596407552675339; F = 1.596407552675339 }
    9  = @{ E = 1.596407552675339; F = 30.2665840088855  }
    10 = @{ E = 12.32790358674824; F = 1.596407552675339 }
    11 = @{ E = 1.596407552675339; F = 1.596407552675339 }
    12 = @{ E = 1.596407552675339; F = 1.596407552675339 }
    13 = @{ E = 1.596407552675339; F = 1.596407552675339 }
    14 = @{ E = 1.596407552675339; F = 1.596407552675339 }
    15 = @{ E = 1.596407552675339; F = 1.596407552675339 }
    16 = @{ E = 1.596407552675339; F = 1.596407552675339 }
    17 = @{ E = 1.596407552675339; F = 1.596407552675339 }
    18 = @{ E = 1.596407552675339; F = 1.596407552675339 }
    19 = @{ E = 1.596407552675339; F = 1.596407552675339 }
    20 = @{ E = 1.596407552675339; F = 1.596407552675339 }
    21 = @{ E = 1.596407552675339; F = 22.75342229809739 }
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row].E
    $ws.Range("F$row").Value = $values[$row].F
}
